$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder
$ws.Range("H6").Value = 166.16667
$ws.Range("I6").Value = 166.16667
$ws.Range("K6").Value = 498.50001
$ws.Range("M6").Value = -386.50001

# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 686.1111
$ws.Range("I38").Value = 146.875
$ws.Range("K38").Value = 440.625
$ws.Range("M38").Value = -68.625

# Row 39: Riches' Brew
$ws.Range("H39").Value = 32.6875
$ws.Range("I39").Value = 15.333333
$ws.Range("K39").Value = 45.999999
$ws.Range("M39").Value = 250.000001

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 11333
$ws.Range("I43").Value = 4333.3335
$ws.Range("K43").Value = 4333.3335
$ws.Range("M43").Value = -4264.3335

# Row 51: A Bile Business
$ws.Range("H51").Value = 46182.36
$ws.Range("J51").Value = 60187.223
$ws.Range("L51").Value = 60187.223
$ws.Range("N51").Value = -61155.223

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 421.2857
$ws.Range("I92").Value = 336.9
$ws.Range("K92").Value = 336.9
$ws.Range("M92").Value = 911.1

# Row 123: Nearly Bare
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3214.5
$ws.Range("I138").Value = 1807.3572
$ws.Range("J138").Value = 6497.8335
$ws.Range("K138").Value = 5422.071599999999
$ws.Range("L138").Value = 19493.5005
$ws.Range("M138").Value = -282.0715999999993
$ws.Range("N138").Value = -29773.5005

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3851.6743
$ws.Range("I32").Value = 2480.6
$ws.Range("K32").Value = 2480.6
$ws.Range("M32").Value = -2193.6

# Row 97: Ore for Me
$ws.Range("H97").Value = 428.05554
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 1000
$ws.Range("N97").Value = -1992

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 863.5
$ws.Range("I22").Value = 227
$ws.Range("K22").Value = 227
$ws.Range("M22").Value = -54

# Row 94: High Steal
$ws.Range("H94").Value = 2232
$ws.Range("I94").Value = 758.75
$ws.Range("J94").Value = 2968.625
$ws.Range("K94").Value = 758.75
$ws.Range("L94").Value = 2968.625
$ws.Range("M94").Value = -307.75
$ws.Range("N94").Value = -3870.625

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move
$ws.Range("H15").Value = 1803.8889
$ws.Range("I15").Value = 3680
$ws.Range("J15").Value = 303
$ws.Range("K15").Value = 3680
$ws.Range("L15").Value = 303
$ws.Range("M15").Value = -3510
$ws.Range("N15").Value = -643

# Row 93: Reeling for Rods
$ws.Range("H93").Value = 11753
$ws.Range("I93").Value = 11753
$ws.Range("K93").Value = 11753
$ws.Range("M93").Value = -9881

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3373.5
$ws.Range("I134").Value = 2748
$ws.Range("K134").Value = 8244
$ws.Range("M134").Value = -5709

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water
$ws.Range("H4").Value = 36167410
$ws.Range("I4").Value = 44847376
$ws.Range("K4").Value = 134542128
$ws.Range("M4").Value = -134542016

# Row 126: Imperial Palate
$ws.Range("H126").Value = 2330
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1368.5927
$ws.Range("J131").Value = 1409.3478
$ws.Range("L131").Value = 4228.0434
$ws.Range("N131").Value = -14308.0434

# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 2297.5
$ws.Range("I134").Value = 2297.5
$ws.Range("K134").Value = 6892.5
$ws.Range("M134").Value = -1822.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4799.6
$ws.Range("I80").Value = 2999
$ws.Range("J80").Value = 5249.75
$ws.Range("K80").Value = 2999
$ws.Range("L80").Value = 5249.75
$ws.Range("M80").Value = -2001
$ws.Range("N80").Value = -7245.75

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4799.6
$ws.Range("I83").Value = 2999
$ws.Range("J83").Value = 5249.75
$ws.Range("K83").Value = 14995
$ws.Range("L83").Value = 26248.75
$ws.Range("M83").Value = -10003
$ws.Range("N83").Value = -36232.75

# Row 94: Wants and Needles
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352

# Row 132: On Board for Lar
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 6465.6665
$ws.Range("I7").Value = 6398.778
$ws.Range("K7").Value = 6398.778
$ws.Range("M7").Value = -6286.778

# Row 18: Simply the Best
$ws.Range("H18").Value = 48000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 23: Back in the Band
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

# Row 41: The Hand that Bleeds
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# Row 42: Slave to Fashion
$ws.Range("H42").Value = 20000
$ws.Range("J42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("N42").Value = -21126

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 598.5
$ws.Range("I46").Value = 495
$ws.Range("J46").Value = 610
$ws.Range("K46").Value = 495
$ws.Range("L46").Value = 610
$ws.Range("M46").Value = -307
$ws.Range("N46").Value = -986

# Row 49: First They Came for the Heretics
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20294

# Row 126: Battered Books
$ws.Range("H126").Value = 6465.6665
$ws.Range("I126").Value = 6398.778
$ws.Range("K126").Value = 19196.334
$ws.Range("M126").Value = -16726.334

$ws = $wb.Worksheets.Item("WVR")
# Row 11: Wiggle Room
$ws.Range("H11").Value = 17996.666
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 17996.666
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 17996.666
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -18280.666

# Row 40: Our Man in Ul'dah
$ws.Range("H40").Value = 10018
$ws.Range("J40").Value = 10018
$ws.Range("L40").Value = 10018
$ws.Range("N40").Value = -10316

# Row 64: Ribbon of Remembrance
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496

# Row 67: The Road Was a Ribbon of Moonlight (L)
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716

# Row 113: A Tender Table
$ws.Range("H113").Value = 2527.0908
$ws.Range("I113").Value = 699.75
$ws.Range("J113").Value = 3571.2856
$ws.Range("K113").Value = 2099.25
$ws.Range("L113").Value = 10713.8568
$ws.Range("M113").Value = 70.75
$ws.Range("N113").Value = -15053.8568

# Row 137: Traditional Trousers
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
